# AHP workbook update: rename criteria labels, refresh the pairwise
# comparison values, and add a new two-column I10:J14 block (with its own
# font/wrap formatting) to support dynamically-loaded participant data -
# matching the rework described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Rename the criteria / "wrt" labels ---
$ws.Range("B3").Value = "Cost"
$ws.Range("B4").Value = "Sustainability"
$ws.Range("B5").Value = "Time to complete mission"
$ws.Range("B6").Value = "Rough Sea Tolerance"
$ws.Range("B7").Value = "Design Complexity"

$ws.Range("C2").Value = "wrt Cost"
$ws.Range("D2").Value = "wrt Sustainability"
$ws.Range("E2").Value = "wrt Time"
$ws.Range("F2").Value = "wrt Rough Sea"
$ws.Range("G2").Value = "wrt Complexity"

# --- 2. Update the pairwise-comparison matrix values ---
$ws.Range("D3").Value = 0.25
$ws.Range("E3").Value = 0.5
$ws.Range("F3").Value = 2

$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 5
$ws.Range("G4").Value = 4

$ws.Range("F5").Value = 3
$ws.Range("G5").Value = 2

# --- 3. Extend the existing number format into column H and row 8 spacer ---
$ws.Range("H2:H7").NumberFormat = "0.00"
$ws.Range("B8").NumberFormat = "0.00"

# --- 4. New I10:J14 block: wrap-text column + plain column, new black font ---
# Build the wrap+font style first (so it is allocated before the plain
# font-only style), then fill in the remaining cells with the plain style.
$ws.Range("J10").Font.Color = 0
$ws.Range("J10").WrapText = $true

$ws.Range("I10").Font.Color = 0
$ws.Range("I10").WrapText = $true

$ws.Range("I11:I14").Font.Color = 0
$ws.Range("J11:J14").Font.Color = 0
$ws.Range("J11:J14").WrapText = $true

# --- 5. Column widths (character widths; Excel snaps these to whole
#        pixels internally, so we pick the inputs that land on the
#        intended ~21.55 / ~13.33 / ~36.89 character widths) ---
$ws.Range("B:B").ColumnWidth = 20.6
$ws.Range("C:G").ColumnWidth = 12.5
$ws.Range("I:I").ColumnWidth = 20.6
$ws.Range("J:J").ColumnWidth = 36.0

# --- 6. Selection to match the saved view state ---
$ws.Range("I10:J14").Select()
